$d = $word.ActiveDocument

$d.Paragraphs.Item(1).Range.Text = "RON LIDGI"
$d.Paragraphs.Item(5).Range.Text = "B.S. Business Administration"
$d.Paragraphs.Item(6).Range.Text = "University of California, Berkeley"
$d.Paragraphs.Item(6).Range.InsertParagraphAfter()
$d.Paragraphs.Item(8).Range.Text = "M.S. Statistics with Actuarial Science Concentration (Online, Part-time)"
$d.Paragraphs.Item(8).Range.InsertParagraphAfter()
$d.Paragraphs.Item(9).Range.Text = "California State University, East Bay"
$d.Paragraphs.Item(9).Range.InsertParagraphAfter()
$d.Paragraphs.Item(10).Range.Text = "Expected Completion: May 2026"
$d.Paragraphs.Item(10).Range.InsertParagraphAfter()
$d.Paragraphs.Item(11).Range.Text = "GPA: 3.6"
$d.Paragraphs.Item(14).Range.Text = "Probability (P)"
$d.Paragraphs.Item(15).Range.Text = "Financial Mathematics (FM)"
$d.Paragraphs.Item(16).Range.Text = "Fundamentals of Actuarial Mathematics (FAM)"
$d.Paragraphs.Item(18).Range.Text = "PROFESSIONAL EXPERIENCE"
$d.Paragraphs.Item(19).Range.Text = "Freelance Data Analyst and Actuarial Consultant"
$d.Paragraphs.Item(19).Range.InsertParagraphAfter()
$d.Paragraphs.Item(20).Range.Text = "Fiverr.com/Upwork.com | July 2024 - Present"
$d.Paragraphs.Item(21).Range.Text = "Delivered statistical analyses and predictive modeling for clients in diverse industries"
$d.Paragraphs.Item(22).Range.Text = "Advised clients on AI-enhanced strategies"
$d.Paragraphs.Item(24).Range.Text = "Intern"
$d.Paragraphs.Item(24).Range.InsertParagraphAfter()
$d.Paragraphs.Item(25).Range.Text = "Casualty Actuarial Society | June 2025 - August 2025"
$d.Paragraphs.Item(26).Range.Text = "Underwent an 8-week professional training program, covering insurance fundamentals"
$d.Paragraphs.Item(26).Range.InsertParagraphAfter()
$d.Paragraphs.Item(27).Range.Text = "Acquired knowledge about property and casualty insurance, Excel, Data Visualization, Ratemaking, Reserving, Predictive Modeling, and essential Soft Skills"
$d.Paragraphs.Item(28).Range.Delete()
$d.Paragraphs.Item(30).Range.Text = "Reinsurance Analysis: Actuarial Case Competition, Traveler’s | Spring 2025"
$d.Paragraphs.Item(31).Range.Text = "Estimated prospective premiums and losses using actuarial methods and historical data"
$d.Paragraphs.Item(32).Range.Text = "Provided strategic recommendations on reinsurance quotes and treaty structures"
$d.Paragraphs.Item(33).Range.Delete()
$d.Paragraphs.Item(33).Range.Delete()
$d.Paragraphs.Item(33).Range.Text = "Conducted Regression Modeling to predict medical expenses"
$d.Paragraphs.Item(34).Range.Text = "Applied EDA, feature selection and data visualization for accurate predictions"
$d.Paragraphs.Item(37).Range.Text = "Python, R, SQL, SAS, Microsoft Excel, Power BI, Tableau, Flask, Web App Development"
$d.Paragraphs.Item(38).Range.Text = "Data Science and Machine Learning: Data Mining, NLP, Keras, TensorFlow, Pytorch, LangChain, Hugging Face"
$d.Paragraphs.Item(39).Range.Text = "Cloud Computing: Azure, AWS Sagemaker, Google Vertex AI"
